$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 5.1
$ws.Range("E2").Value = 0.09
$ws.Range("F2").Value = 307
$ws.Range("G2").Value = 0.0001
$ws.Range("I2").Value = 1000
$ws.Range("N2").Value = "[0.02]"
